$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.228.35'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.029.62'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.63'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.31'
$ws.Range("E6").Value = '  +3.24%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.023.65'
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.67'
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.480'
$ws.Range("E12").Value = '  +5.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("E13").Value = '  -2.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.38'
$ws.Range("E14").Value = '  +5.04%  '
$ws.Range("E15").Value = '  -0.30%  '
$ws.Range("D16").Value = '66.206.85'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = '3.533.41'
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.23'
$ws.Range("E18").Value = '  +4.64%  '
$ws.Range("E19").Value = '  +18.79%  '
$ws.Range("D20").Value = '3.032.83'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '471.46'
$ws.Range("E21").Value = '  +3.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.706'
$ws.Range("E23").Value = '  +1.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.25'
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.87'
$ws.Range("E25").Value = '  +5.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.05'
$ws.Range("E27").Value = '  -4.02%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.21'
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.45'
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.63'
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.119'
$ws.Range("E32").Value = '  +7.68%  '
$ws.Range("D33").Value = '0.0₃0996'
$ws.Range("E33").Value = '  -6.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.99'
$ws.Range("E34").Value = '  +2.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.992'
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.85'
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.20'
$ws.Range("E38").Value = '  +10.27%  '
$ws.Range("E39").Value = '  -5.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.65'
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("E42").Value = '  -0.79%  '
$ws.Range("E43").Value = '  -5.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.60'
$ws.Range("E44").Value = '  +2.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0360'
$ws.Range("E45").Value = '  -0.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '382.41'
$ws.Range("E46").Value = '  -4.08%  '
$ws.Range("D47").Value = '2.722.46'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.38'
$ws.Range("E48").Value = '  +0.07%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.56'
$ws.Range("E50").Value = '  +2.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.22'
$ws.Range("E51").Value = '  +3.94%  '
